# Journal-Entry-template.xlsx: update COA Code column (D) amounts and
# move the saved cell selection from D10 to D5 (matches the author's
# final on-screen selection before saving).
#
# NOTE: the commit also rewrote several purely cosmetic / environment
# artifacts that Excel stamps into the OOXML on save from this machine
# (the mc:Choice x15ac:absPath folder, the xr:revisionPtr documentId
# GUID, and the bookViews window position/size) and reshuffled the
# customXml SharePoint metadata parts (item2.xml <-> item3.xml). None of
# those are reachable through the Excel object model -- they are not
# values a user (or macro) sets via Application/Workbook/Worksheet/Range
# properties, so there is nothing equivalent to drive here; they are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# COA Code column updates (D2:D5)
$ws.Range("D2").Value = 190
$ws.Range("D3").Value = 501
$ws.Range("D4").Value = 160
$ws.Range("D5").Value = 314

# Move the active selection from D10 to D5
[void]$ws.Range("D5").Select()
